# REVER_DailyTracker_MONISHA.xlsx -- "Add files via upload"
#
# The "JUNE-2021" tab actually held May-2021 dates (44317 = 2021-05-01 ..
# 44347 = 2021-05-31, 31 rows). The edit shifts every date by +31 days so
# the tab finally holds real June-2021 dates (44348 = 2021-06-01 ..
# 44377 = 2021-06-30, 30 rows), fills in the day that is now the last row
# of the month (row 31 / June 30) with real task data, and blanks out the
# now-unused 31st slot (row 32, since June only has 30 days).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JUNE-2021")

# --- Shift the daily date column (B2:B30) forward by one month (31 days) ---
# May 1 (44317) -> June 1 (44348), ... May 29 (44345) -> June 29 (44376).
# (Written as literal target serials -- this host's Range.Value getter does
# not reliably round-trip, so read-modify-write is avoided here.)
$juneDates = @{
    2 = 44348; 3 = 44349; 4 = 44350; 5 = 44351; 6 = 44352
    7 = 44353; 8 = 44354; 9 = 44355; 10 = 44356; 11 = 44357
    12 = 44358; 13 = 44359; 14 = 44360; 15 = 44361; 16 = 44362
    17 = 44363; 18 = 44364; 19 = 44365; 20 = 44366; 21 = 44367
    22 = 44368; 23 = 44369; 24 = 44370; 25 = 44371; 26 = 44372
    27 = 44373; 28 = 44374; 29 = 44375; 30 = 44376
}
foreach ($row in $juneDates.Keys) {
    $ws.Cells.Item($row, 2).Value = $juneDates[$row]
}

# --- Row 31: was the blank "May 31" placeholder row, now becomes the
#     real "June 30" entry (day 30, task done, status Completed) ---
$ws.Range("B31").Value = 44377

# Pull formatting for the newly-populated cells from row 30, which already
# carries the "filled task row" styles we need (s=13 / s=50 / s=54).
$ws.Range("C30").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D30").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("F26").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C31").Value = "Hayaai"
$ws.Range("D31").Value = "Inventory Sceens modifications done"
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Completed"

# --- Row 32: June only has 30 days, so the old "day 31" row (A32/B32)
#     is cleared out, leaving the rest of the row untouched ---
$ws.Range("A32:B32").ClearContents()

# --- Restore the sheet selection to match where the editor left off ---
$ws.Range("D31").Select()
